$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 5 per the task schedule change
$ws.Range("D5").Value = 2
$ws.Range("F5").Value = -3
$ws.Range("H5").Value = 46

# Move the active selection to D5
$ws.Range("D5").Select()
